# Update column F (dSF) values on the active sheet to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = -1
    4  = 3
    5  = 3
    6  = 1
    7  = 4
    8  = 1
    10 = 2
    11 = -1
    12 = 5
    13 = 1
    14 = 2
    16 = -1
    17 = -4
    18 = 4
    20 = 5
    21 = 5
    22 = 1
    23 = -5
    24 = -3
    25 = -1
    26 = 1
    27 = -1
    28 = -1
    30 = 4
    31 = 3
    32 = -3
    33 = 2
    34 = -1
    35 = -1
    36 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$wb.Save()
